$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the "Save" header in H1, matching style of existing header cells (B1:G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the new "Save" column values for rows 2-7
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 0
